# Scheduled market-data refresh: update Market Board price/profit columns
# (H:currentAveragePrice, I:currentAveragePriceNQ, J:currentAveragePriceHQ,
#  K:LevePriceNQ, L:LevePriceHQ, M:LeveProfitNQ, N:LeveProfitHQ) per sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2314.1428
$ws.Range("J28").Value = 3166.3333
$ws.Range("L28").Value = 3166.3333
$ws.Range("N28").Value = -4136.3333

$ws.Range("H69").Value = 3500
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 3500
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

$ws.Range("H106").Value = 4995
$ws.Range("I106").Value = 4995
$ws.Range("K106").Value = 4995
$ws.Range("M106").Value = -4364

$ws.Range("H138").Value = 4770.227
$ws.Range("I138").Value = 4197.7144
$ws.Range("J138").Value = 5037.4
$ws.Range("K138").Value = 12593.1432
$ws.Range("L138").Value = 15112.2
$ws.Range("M138").Value = -7453.143199999999
$ws.Range("N138").Value = -25392.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2424.7144
$ws.Range("I2").Value = 2242.5
$ws.Range("K2").Value = 2242.5
$ws.Range("M2").Value = -2129.5

$ws.Range("H32").Value = 13560.958
$ws.Range("I32").Value = 12903.381
$ws.Range("J32").Value = 18164
$ws.Range("K32").Value = 12903.381
$ws.Range("L32").Value = 18164
$ws.Range("M32").Value = -12616.381
$ws.Range("N32").Value = -18738

$ws.Range("H61").Value = 1732.6
$ws.Range("I61").Value = 1319.8948
$ws.Range("K61").Value = 1319.8948
$ws.Range("M61").Value = -1107.8948

$ws.Range("H116").Value = 2424.7144
$ws.Range("I116").Value = 2242.5
$ws.Range("K116").Value = 2242.5
$ws.Range("M116").Value = 51.5

$ws.Range("H132").Value = 2094.0557
$ws.Range("I132").Value = 1446.3334
$ws.Range("K132").Value = 4339.0002
$ws.Range("M132").Value = -1809.0002

$ws.Range("H136").Value = 1732.6
$ws.Range("I136").Value = 1319.8948
$ws.Range("K136").Value = 3959.6844
$ws.Range("M136").Value = -1409.6844

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2424.7144
$ws.Range("I3").Value = 2242.5
$ws.Range("K3").Value = 2242.5
$ws.Range("M3").Value = -2128.5

$ws.Range("H86").Value = 2074.9524
$ws.Range("I86").Value = 2276.1765
$ws.Range("J86").Value = 1219.75
$ws.Range("K86").Value = 2276.1765
$ws.Range("L86").Value = 1219.75
$ws.Range("M86").Value = -1153.1765
$ws.Range("N86").Value = -3465.75

$ws.Range("H89").Value = 2074.9524
$ws.Range("I89").Value = 2276.1765
$ws.Range("J89").Value = 1219.75
$ws.Range("K89").Value = 11380.8825
$ws.Range("L89").Value = 6098.75
$ws.Range("M89").Value = -5764.8825
$ws.Range("N89").Value = -17330.75

$ws.Range("H134").Value = 3128.8572
$ws.Range("I134").Value = 2981.875
$ws.Range("K134").Value = 8945.625
$ws.Range("M134").Value = -6410.625

$ws.Range("H140").Value = 94966.664
$ws.Range("J140").Value = 94966.664
$ws.Range("L140").Value = 94966.664
$ws.Range("N140").Value = -105326.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 626.6667
$ws.Range("J22").Value = 580
$ws.Range("L22").Value = 580
$ws.Range("N22").Value = -1280

$ws.Range("H31").Value = 2171.5217
$ws.Range("I31").Value = 2164.1428
$ws.Range("J31").Value = 2249
$ws.Range("K31").Value = 2164.1428
$ws.Range("L31").Value = 2249
$ws.Range("M31").Value = -1869.1428
$ws.Range("N31").Value = -2839

$ws.Range("H34").Value = 2171.5217
$ws.Range("I34").Value = 2164.1428
$ws.Range("J34").Value = 2249
$ws.Range("K34").Value = 2164.1428
$ws.Range("L34").Value = 2249
$ws.Range("M34").Value = -1962.1428
$ws.Range("N34").Value = -2653

$ws.Range("H58").Value = 2440.6924
$ws.Range("I58").Value = 2339.0908
$ws.Range("K58").Value = 2339.0908
$ws.Range("M58").Value = -2136.0908

$ws.Range("H86").Value = 23077.842
$ws.Range("J86").Value = 36947.555
$ws.Range("L86").Value = 36947.555
$ws.Range("N86").Value = -39193.555

$ws.Range("H89").Value = 23077.842
$ws.Range("J89").Value = 36947.555
$ws.Range("L89").Value = 184737.775
$ws.Range("N89").Value = -195969.775

$ws.Range("H107").Value = 1848.7693
$ws.Range("I107").Value = 1406.3334
$ws.Range("J107").Value = 2228
$ws.Range("K107").Value = 1406.3334
$ws.Range("L107").Value = 2228
$ws.Range("M107").Value = 513.6666
$ws.Range("N107").Value = -6068

$ws.Range("H132").Value = 4762.5
$ws.Range("I132").Value = 4728.7144
$ws.Range("K132").Value = 14186.1432
$ws.Range("M132").Value = -11656.1432

$ws.Range("H136").Value = 2440.6924
$ws.Range("I136").Value = 2339.0908
$ws.Range("K136").Value = 7017.2724
$ws.Range("M136").Value = -4467.2724

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3621.875
$ws.Range("J113").Value = 3568.5715
$ws.Range("L113").Value = 10705.7145
$ws.Range("N113").Value = -15045.7145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 91.75
$ws.Range("I2").Value = 101.28571
$ws.Range("J2").Value = 25
$ws.Range("K2").Value = 101.28571
$ws.Range("L2").Value = 25
$ws.Range("M2").Value = 11.71429000000001
$ws.Range("N2").Value = -251

$ws.Range("H14").Value = 97.5
$ws.Range("I14").Value = 95
$ws.Range("J14").Value = 100
$ws.Range("K14").Value = 95
$ws.Range("L14").Value = 100
$ws.Range("M14").Value = 73
$ws.Range("N14").Value = -436

$ws.Range("H43").Value = 3273.8
$ws.Range("J43").Value = 3842.25
$ws.Range("L43").Value = 3842.25
$ws.Range("N43").Value = -4144.25

$ws.Range("H46").Value = 4124.2856
$ws.Range("J46").Value = 4124.2856
$ws.Range("L46").Value = 4124.2856
$ws.Range("N46").Value = -4436.2856

$ws.Range("H57").Value = 13329
$ws.Range("J57").Value = 16497.5
$ws.Range("L57").Value = 16497.5
$ws.Range("N57").Value = -18137.5

$ws.Range("H80").Value = 6683.857
$ws.Range("J80").Value = 7058
$ws.Range("L80").Value = 7058
$ws.Range("N80").Value = -9054

$ws.Range("H83").Value = 6683.857
$ws.Range("J83").Value = 7058
$ws.Range("L83").Value = 35290
$ws.Range("N83").Value = -45274

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6191.778
$ws.Range("J22").Value = 6789
$ws.Range("L22").Value = 6789
$ws.Range("N22").Value = -7379

$ws.Range("H27").Value = 6191.778
$ws.Range("J27").Value = 6789
$ws.Range("L27").Value = 6789
$ws.Range("N27").Value = -7003

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1817.3077
$ws.Range("I136").Value = 1512.7
$ws.Range("J136").Value = 2832.6667
$ws.Range("K136").Value = 4538.1
$ws.Range("L136").Value = 8498.000100000001
$ws.Range("M136").Value = -1988.1
$ws.Range("N136").Value = -13598.0001
